# Re-order a handful of item/UOM rows in the NoStock sheet.
# (Matches author commit: items/UOM re-sequenced for specific brand groups;
#  only the Item Name (D) / UOM (E) cell values for the affected rows move -
#  everything else on the sheet is untouched.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Dinafex 180mg Tablet"
$ws.Range("D5").Value = "Dinafex 60mg Tablet"
$ws.Range("D7").Value = "Etorix 90mg Tablet"
$ws.Range("E7").Value = "30's"
$ws.Range("D9").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E9").Value = "40's"
$ws.Range("D11").Value = "Flucloxin 500mg Capsule"
$ws.Range("E11").Value = "30 's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E12").Value = "36 's"
$ws.Range("D14").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E14").Value = "4's"
$ws.Range("D15").Value = "Ketonic 10mg Tablet"
$ws.Range("E15").Value = "20's"
$ws.Range("D16").Value = "Ketonic 30mg Injection"
$ws.Range("E16").Value = "5 's"
$ws.Range("D18").Value = "Kynol TR 100mg Capsule"
$ws.Range("E18").Value = "50 's"
$ws.Range("D19").Value = "Kynol TR 200mg Capsule"
$ws.Range("E19").Value = "30 's"
$ws.Range("D24").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E24").Value = "30ml"
$ws.Range("D25").Value = "Zithrox 15ml Suspension"
$ws.Range("E25").Value = "15 ml"
$ws.Range("D26").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E26").Value = "6's"
$ws.Range("D27").Value = "Zithrox 500mg Tablet"
$ws.Range("E27").Value = "6 's"
